$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.985.93'
$ws.Range("E2").Value = '  -2.99%  '

$ws.Range("D3").Value = '1.890.67'
$ws.Range("E3").Value = '  -3.84%  '

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.005'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.66%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '326.36'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.70%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '1.005'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.63%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.4585'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -3.97%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3935'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -2.81%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '51.50'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -4.66%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.08220'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -3.79%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '1.035'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -2.55%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '21.62'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -3.90%  '

$ws.Range("D13").Value = '1.875.26'
$ws.Range("E13").Value = '  -4.53%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '7.309'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -4.77%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '5.969'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -4.66%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '1.007'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.77%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '89.10'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -1.17%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.00001059'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -1.14%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.06567'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -0.81%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '17.62'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -5.80%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -1.03%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '5.651'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -2.46%  '

$ws.Range("D23").Value = '28.001.56'
$ws.Range("E23").Value = '  -2.94%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '11.08'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -4.42%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.308'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.68%  '

$ws.Range("D26").Value = '2.141.53'
$ws.Range("E26").Value = '  -2.53%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '153.96'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.21%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '19.91'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -1.75%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '2.099'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -2.48%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '5.654'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -5.35%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '123.97'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -0.40%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.09533'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -1.10%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.9584'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -5.04%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.455'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -0.73%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '3.631'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -1.60%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '5.454'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -4.18%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.02288'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -2.58%  '

$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '1.255'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -1.16%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '8.659'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -1.33%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.06103'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -1.83%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.6092'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -2.68%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '1.004'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -0.66%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '10.70'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -3.73%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.1887'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -1.99%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '1.308'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -3.21%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.5811'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -2.93%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '12.67'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -2.49%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '1.990'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -4.75%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '3.427'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +0.08%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.06889'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +0.45%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '110.33'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.72%  '
